$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Hours=1, Start time="11pm", End time="11:54pm"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "11pm"
$ws.Range("D8").Value = "11:54pm"

# Row 9: Date=1/28/2020, Hours=2, Start time="1:45am", End time="3:45am"
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 43858
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "1:45am"
$ws.Range("D9").Value = "3:45am"

$ws.Range("D13").Select()
